$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Output")

# Row 2 - update timestamp only (rest unchanged)
$ws.Range("A2").Value = "21/01/2022 12:00:21 pm"

# Row 3 - update timestamp only (rest unchanged); style on message cell reverts to default
$ws.Range("A3").Value = "21/01/2022 12:00:28 pm"
$ws.Range("C3").Style = "Normal"

# Row 4 - now an error row (Fare elements not found. Timeout); clear fare columns
$ws.Range("A4").Value = "21/01/2022 12:01:01 pm"
$ws.Range("B4").Value = "Y"
$ws.Range("C4").Value = "Fare elements not found. Timeout"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4:G4").ClearContents()

# Row 5 - update timestamp only (rest unchanged); style on message cell reverts to default
$ws.Range("A5").Value = "21/01/2022 12:01:09 pm"
$ws.Range("C5").Style = "Normal"

# Row 6 - update timestamp and fare amounts; style on message cell reverts to default
$ws.Range("A6").Value = "21/01/2022 12:01:22 pm"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "₹1,44,570"
$ws.Range("E6").Value = "₹10,340"
$ws.Range("G6").Value = "₹1,54,920"

# Row 7 - now an error row (Flight type not found); clear fare columns
$ws.Range("A7").Value = "21/01/2022 11:45:51 am"
$ws.Range("B7").Value = "Y"
$ws.Range("C7").Value = "Flight type not found"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7:G7").ClearContents()
